$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels: C1 becomes "Resolución Primigenia",
# D1 becomes "Resolución Hija"
$ws.Range("C1").Value = "Resolución Primigenia"
$ws.Range("D1").Value = "Resolución Hija"

# Clear out the "Resolución Hija" values in D2 and D3 (now blank,
# since they previously duplicated the parent resolution number)
$ws.Range("D2:D3").ClearContents()
